$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") contains a date value that was bumped from
# 45192 (2023-09-23) to 45202 (2023-10-03) for every data row (2-319).
$ws.Range("C2:C319").Value = 45202
